# Apply the "no is_pref / no lev distance" re-export:
#  - column B (id) / C (speaker_variant) values are re-ordered across rows 2-38
#  - column D (is_prefered) is cleared for every data row (no more "x" markers)
#  - a new row 38 is appended (duplicate of the old row 26 entry)
#  - sheet dimension grows from A1:H37 to A1:H38

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$url = "https://www.dbnl.org/tekst/nieu001jeru01_01"

# New id (col B) / speaker_variant (col C) values, row 2 .. row 38
$data = @(
    @{row=2;  id="#iuchal";           variant="Iuchal"}
    @{row=3;  id="#c.-zedechias";     variant="C. Zedechias"}
    @{row=4;  id="#nabuc";            variant="Nabuc"}
    @{row=5;  id="#nabusa";           variant="Nabusa"}
    @{row=6;  id="#nabusardan";       variant="Nabusardan"}
    @{row=7;  id="#g.-gods";          variant="G. Gods"}
    @{row=8;  id="#melinda";          variant="Melinda"}
    @{row=9;  id="#gedalia";          variant="Gedalia"}
    @{row=10; id="#hamiltal";         variant="Hamiltal"}
    @{row=11; id="#1.-sone";          variant="1. Sone"}
    @{row=12; id="#g.-godts";         variant="G. Godts"}
    @{row=13; id="#binnen-nabuch";    variant="Binnen Nabuch"}
    @{row=14; id="#ieremias";         variant="Ieremias"}
    @{row=15; id="#thessa";           variant="Thessa"}
    @{row=16; id="#1.-soone";         variant="1. Soone"}
    @{row=17; id="#wt";               variant="Wt"}
    @{row=18; id="#ebedme";           variant="Ebedme"}
    @{row=19; id="#zedech";           variant="Zedech"}
    @{row=20; id="#zarezar";          variant="Zarezar"}
    @{row=21; id="#nabuchodonosor";   variant="Nabuchodonosor"}
    @{row=22; id="#zarezar-wt";       variant="Zarezar wt"}
    @{row=23; id="#zedechias";        variant="Zedechias"}
    @{row=24; id="#nebo";             variant="Nebo"}
    @{row=25; id="#1.-soon";          variant="1. Soon"}
    @{row=26; id="#nargal";           variant="Nargal"}
    @{row=27; id="#2.-soon";          variant="2. Soon"}
    @{row=28; id="#nabuch";           variant="Nabuch"}
    @{row=29; id="#nabusar";          variant="Nabusar"}
    @{row=30; id="#ieremias-wt";      variant="Ieremias wt"}
    @{row=31; id="#thessalia";        variant="Thessalia"}
    @{row=32; id="#ieremias-binnen";  variant="Ieremias binnen"}
    @{row=33; id="#gadalia";          variant="Gadalia"}
    @{row=34; id="#thessal";          variant="Thessal"}
    @{row=35; id="#pashur";           variant="Pashur"}
    @{row=36; id="#2";                variant="2"}
    @{row=37; id="#2-sone";           variant="2 Sone"}
    @{row=38; id="#nabusar.-binnen";  variant="Nabusar. binnen"}
)

foreach ($item in $data) {
    $r = $item.row
    $ws.Cells.Item($r, 1).Value2 = $url
    $ws.Cells.Item($r, 2).Value2 = $item.id
    $variant = $item.variant
    if ($variant -match '^[0-9]+(\.[0-9]+)?$') {
        # Force text so a purely-numeric-looking label (e.g. "2") is not
        # silently turned into a number, same as typing a leading apostrophe
        # would do in the Excel UI.
        $ws.Cells.Item($r, 3).Value2 = "'" + $variant
    } else {
        $ws.Cells.Item($r, 3).Value2 = $variant
    }
}

# Column D ("is_prefered") no longer carries any "x" markers - clear it for
# every data row (rows 2-29 previously held "x").
$ws.Range("D2:D29").ClearContents()

Write-Host "done"
